$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the existing "1:14" entry (row 20) - it was missing "AM"
$ws.Range("B20").Value = "1:14AM 1-17-2018"

# Turn the old totals row (row 21) into a new data row
$ws.Range("A21").Value = "6:33PM 1-19-2018"
$ws.Range("B21").Value = "9:06PM 1-19-2018"
$ws.Range("C21").Value = 153

# Move "Total Project Hours:" down to row 29 with an updated SUM formula
$ws.Range("A29").Value = "Total Project Hours:"
$ws.Range("C29").Formula = "=SUM(C2:C28)/60"

# Scroll the view down a bit and leave the selection on C22, matching where
# editing left off
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C22").Select()
